# "updated user story test cases"
#
# Appends a new "User Story" block of test cases (TC-083 .. TC-102) to the
# bottom of the existing Test Id / Test Scenario table on Sheet1, growing
# the used range from A1:G82 to A1:G102.
#
# The cell-by-cell order below intentionally mirrors how the rows were
# originally authored (column B's single header cell, then column C top
# to bottom, then column A top to bottom, with a couple of out-of-order
# tail entries) so the workbook's shared-string table ends up in the same
# append order as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A ("Test Id") for rows 2-82 uses a bordered cell style that isn't
# the worksheet's default. Copy that formatting from the last existing row
# down across the new rows before typing values, so A83:A102 pick up the
# same style as A82 instead of the unstyled default.
$ws.Cells.Item(82, 1).Copy() | Out-Null
$ws.Range("A83:A102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New section header in column B.
$ws.Cells.Item(83, 2).Value = "User Story"

# Column C ("Test Scenario") filled top-down for rows 83-99.
$ws.Cells.Item(83, 3).Value = "Login page to register user should be available"
$ws.Cells.Item(84, 3).Value = "Login page to logged in as registered user"
$ws.Cells.Item(85, 3).Value = "On register page option to provide name, username, password and email id is available"
$ws.Cells.Item(86, 3).Value = "Verify name fields allows upto 50 characters"
$ws.Cells.Item(87, 3).Value = "Verify username fields allows upto 50 characters"
$ws.Cells.Item(88, 3).Value = "Verify password is minimum of 15 characters"
$ws.Cells.Item(89, 3).Value = "Verfiy password should not start with special characters"
$ws.Cells.Item(90, 3).Value = "Verify  Email id is in format XXX@XXX.com"
$ws.Cells.Item(91, 3).Value = "Login with the user created above and create the workflows"
$ws.Cells.Item(92, 3).Value = "Viewing of the workflow is possible"
$ws.Cells.Item(93, 3).Value = "Adding workflow is possible"
$ws.Cells.Item(94, 3).Value = "Modifying workflow is possible"
$ws.Cells.Item(95, 3).Value = "Deleting the workflow is possible"
$ws.Cells.Item(96, 3).Value = "Logged in user name should be displayed on the worklfow screen"
$ws.Cells.Item(97, 3).Value = "Logged in user name should be displayed on the mested worklfow screen"
$ws.Cells.Item(98, 3).Value = "Create another user and loggin with it. Workflows created by 1st user should not be visible to the 2nd user"
$ws.Cells.Item(99, 3).Value = "Logout button should be available"

# Column A ("Test Id") filled top-down for rows 83-100.
$ws.Cells.Item(83, 1).Value = "TC-083"
$ws.Cells.Item(84, 1).Value = "TC-084"
$ws.Cells.Item(85, 1).Value = "TC-085"
$ws.Cells.Item(86, 1).Value = "TC-086"
$ws.Cells.Item(87, 1).Value = "TC-087"
$ws.Cells.Item(88, 1).Value = "TC-088"
$ws.Cells.Item(89, 1).Value = "TC-089"
$ws.Cells.Item(90, 1).Value = "TC-090"
$ws.Cells.Item(91, 1).Value = "TC-091"
$ws.Cells.Item(92, 1).Value = "TC-092"
$ws.Cells.Item(93, 1).Value = "TC-093"
$ws.Cells.Item(94, 1).Value = "TC-094"
$ws.Cells.Item(95, 1).Value = "TC-095"
$ws.Cells.Item(96, 1).Value = "TC-096"
$ws.Cells.Item(97, 1).Value = "TC-097"
$ws.Cells.Item(98, 1).Value = "TC-098"
$ws.Cells.Item(99, 1).Value = "TC-099"
$ws.Cells.Item(100, 1).Value = "TC-100"

# Remaining cells, typed in authoring order: C100, C102, then C101 (which
# reuses the pre-existing "duplicate workflow name" scenario text), and
# finally A101 / A102.
$ws.Cells.Item(100, 3).Value = "On click of the logout user should be navigated back to the login screen"
$ws.Cells.Item(102, 3).Value = "User with same username should not be allowed to regsiter"
$ws.Cells.Item(101, 3).Value = "User should not be allowed to create a new workflow with an existing name"
$ws.Cells.Item(101, 1).Value = "TC-101"
$ws.Cells.Item(102, 1).Value = "TC-102"

# Rows whose Test Scenario text wraps onto a second line get the taller
# 29pt row height, matching the pattern already used elsewhere in the sheet
# (e.g. rows 6 and 41).
$ws.Rows.Item(85).RowHeight = 29
$ws.Rows.Item(97).RowHeight = 29
$ws.Rows.Item(98).RowHeight = 29
$ws.Rows.Item(100).RowHeight = 29
$ws.Rows.Item(101).RowHeight = 29

# Leave the selection on the newly added tail rows, matching where the
# edit left the cursor.
$ws.Range("A100:A102").Select()
